$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.308.08"
$ws.Range("E2").Value = "  -3.19%  "
$ws.Range("D3").Value = "1.931.49"
$ws.Range("E3").Value = "  -3.48%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'249.04"
$ws.Range("E5").Value = "  -3.96%  "
$ws.Range("D6").Value = "'0.7290"
$ws.Range("E6").Value = "  -9.40%  "
$ws.Range("D7").Value = "'0.9994"
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").Value = "'0.3289"
$ws.Range("D9").Value = "'27.38"
$ws.Range("E9").Value = "  -4.15%  "
$ws.Range("D10").Value = "'0.06876"
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("D11").Value = "'0.8080"
$ws.Range("E11").Value = "  -4.26%  "
$ws.Range("D12").Value = "'0.08054"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "1.931.80"
$ws.Range("E13").Value = "  -3.50%  "
$ws.Range("D14").Value = "'5.412"
$ws.Range("E14").Value = "  -4.03%  "
$ws.Range("D15").Value = "'94.99"
$ws.Range("E15").Value = "  -5.83%  "
$ws.Range("D16").Value = "'14.52"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "30.320.11"
$ws.Range("E17").Value = "  -3.13%  "
$ws.Range("D18").Value = "'253.35"
$ws.Range("E18").Value = "  -7.46%  "
$ws.Range("D19").Value = "'0.000008034"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").Value = "'5.829"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("D21").Value = "2.188.50"
$ws.Range("E21").Value = "  -3.28%  "
$ws.Range("D22").Value = "'0.9993"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").Value = "'0.9987"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "'6.869"
$ws.Range("E24").Value = "  -4.38%  "
$ws.Range("D25").Value = "'9.691"
$ws.Range("E25").Value = "  -5.17%  "
$ws.Range("D26").Value = "'159.43"
$ws.Range("E26").Value = "  -2.77%  "
$ws.Range("D27").Value = "'2.397"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "'0.1346"
$ws.Range("E28").Value = "  -13.36%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'19.13"
$ws.Range("E29").Value = "  -4.65%  "
$ws.Range("D30").Value = "'1.557"
$ws.Range("E30").Value = "  -4.46%  "
$ws.Range("D31").Value = "'1.336"
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").Value = "'4.402"
$ws.Range("E32").Value = "  -5.02%  "
$ws.Range("D33").Value = "'4.185"
$ws.Range("E33").Value = "  -4.99%  "
$ws.Range("D34").Value = "'0.05100"
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("D35").Value = "'1.222"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").Value = "'0.7399"
$ws.Range("E36").Value = "  -3.42%  "
$ws.Range("D37").Value = "'2.750"
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("D38").Value = "'0.01969"
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("D39").Value = "'2.820"
$ws.Range("E39").Value = "  -4.28%  "
$ws.Range("D40").Value = "'6.613"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("D41").Value = "'79.24"
$ws.Range("E41").Value = "  -2.07%  "
$ws.Range("D42").Value = "'0.4472"
$ws.Range("E42").Value = "  -5.70%  "
$ws.Range("D43").Value = "'1.997"
$ws.Range("E43").Value = "  -9.71%  "
$ws.Range("D44").Value = "'0.9993"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("E45").Value = "  -3.06%  "
$ws.Range("D46").Value = "'102.18"
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("D47").Value = "'9.793"
$ws.Range("E47").Value = "  -2.01%  "
$ws.Range("D48").Value = "'7.294"
$ws.Range("E48").Value = "  -4.75%  "
$ws.Range("D49").Value = "'36.48"
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").Value = "'0.05944"
$ws.Range("E51").Value = "  -0.33%  "
